$d = $word.ActiveDocument

# Old/new text built with explicit Unicode chars (apostrophe stays straight,
# i.e. we avoid Find/Replace's "replace text" mode which triggers Word's
# smart-quote AutoCorrect and turns ' into a curly quote).
$old = "Dates de la campanya Constel" + [char]0x00B7 + "laci" + [char]0x00F3 + " de Cygnus 2022: 10-19 d'agost, 9-18 de setembre, 8-17 d'octubre"
$new = "Dates de la campanya 2022 en qu" + [char]0x00E8 + " usem la constel" + [char]0x00B7 + "laci" + [char]0x00F3 + ", Constel" + [char]0x00B7 + "laci" + [char]0x00F3 + " de Cygnus 10-19 d'agost, 9-18 de setembre, 8-17 d'octubre"

$r = $d.Content
$r.Find.ClearFormatting()
[void]$r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($r.Find.Found) {
    $r.Text = $new
    $r.Collapse(0)
    $r.End = $d.Content.End
    [void]$r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}
